$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Edit 1: "Client to client authentication" paragraph.
# Split the run "each other. They will then verify them. If both are valid,
# they will use " into several runs so that " establish TLS connection and"
# is woven in (with the _GoBack bookmark sitting right after "establish"),
# matching how Word would have recorded a live edit at that cursor spot.
# ---------------------------------------------------------------------------

$rng = $d.Content
$rng.Find.ClearFormatting()
$found = $rng.Find.Execute("each other. They will then verify them. If both are valid, they will use ")
if (-not $found) {
    throw "Could not locate the 'Client to client authentication' sentence to edit."
}

# Grab the whole enclosing paragraph so the XML splice lands exactly in
# place (InsertXML only behaves predictably over a full paragraph range).
$para = $rng.Paragraphs(1)
$target = $d.Range($para.Range.Start, $para.Range.End)

$xml = @'
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
  <pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
    <pkg:xmlData>
      <w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
        <w:body>
          <w:p>
            <w:r><w:t>Client</w:t></w:r>
            <w:r w:rsidR="00EA5B1C"><w:t>s</w:t></w:r>
            <w:r><w:t xml:space="preserve"> establishing connection will first send its signed </w:t></w:r>
            <w:r w:rsidR="00EA5B1C"><w:t>certificates</w:t></w:r>
            <w:r><w:t xml:space="preserve"> to </w:t></w:r>
            <w:r><w:t>each other. They will then verify them. If both are valid, they will</w:t></w:r>
            <w:r><w:t xml:space="preserve"> </w:t></w:r>
            <w:r><w:t>establish</w:t></w:r>
            <w:bookmarkStart w:id="0" w:name="_GoBack"/>
            <w:bookmarkEnd w:id="0"/>
            <w:r><w:t xml:space="preserve"> TLS connection and</w:t></w:r>
            <w:r><w:t xml:space="preserve"> use </w:t></w:r>
            <w:proofErr w:type="spellStart"/>
            <w:r w:rsidR="00EA5B1C" w:rsidRPr="00EA5B1C"><w:t>Diffie</w:t></w:r>
            <w:proofErr w:type="spellEnd"/>
            <w:r w:rsidR="00EA5B1C" w:rsidRPr="00EA5B1C"><w:t>-Hellman</w:t></w:r>
            <w:r w:rsidR="00EA5B1C"><w:t xml:space="preserve"> key exchange to create encryption and MAC keys. The exchange messages will be</w:t></w:r>
            <w:r w:rsidR="00D12EE5"><w:t xml:space="preserve"> appended with usernames of both clients and timestamp, then</w:t></w:r>
            <w:r w:rsidR="00EA5B1C"><w:t xml:space="preserve"> signed by the </w:t></w:r>
            <w:r w:rsidR="004200D0"><w:t>client’s</w:t></w:r>
            <w:r w:rsidR="00EA5B1C"><w:t xml:space="preserve"> </w:t></w:r>
            <w:r w:rsidR="00D12EE5"><w:t>private key</w:t></w:r>
            <w:r w:rsidR="00EA5B1C"><w:t>. All further messages will</w:t></w:r>
            <w:r w:rsidR="00CA6EC1"><w:t xml:space="preserve"> be encrypted with</w:t></w:r>
            <w:r w:rsidR="00EA5B1C"><w:t xml:space="preserve"> AES-256 and </w:t></w:r>
            <w:r w:rsidR="00CA6EC1"><w:t>authenticated with</w:t></w:r>
            <w:r w:rsidR="00F727B0"><w:t xml:space="preserve"> HMAC-SHA</w:t></w:r>
            <w:r w:rsidR="00AF6AF7"><w:t>512</w:t></w:r>
            <w:r w:rsidR="00F727B0"><w:t xml:space="preserve"> using the newly generated shared secret.</w:t></w:r>
          </w:p>
        </w:body>
      </w:document>
    </pkg:xmlData>
  </pkg:part>
</pkg:package>
'@

$target.InsertXML($xml)

# ---------------------------------------------------------------------------
# Edit 2: the _GoBack bookmark used to live alone in the trailing empty
# paragraph; now that it has moved into the paragraph above, drop the old
# one so the final paragraph goes back to being completely empty.
# ---------------------------------------------------------------------------

try {
    $oldBookmark = $d.Bookmarks("_GoBack")
    $oldBookmark.Delete()
} catch {
    # already gone - nothing to do
}

Write-Host "Edit applied."
